$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00151400454201363
$ws.Range("C2").Value = 0.965177895533687
$ws.Range("D2").Value = 0.00757002271006813
$ws.Range("E2").Value = 0.975018925056775
$ws.Range("F2").Value = 0.00227100681302044
$ws.Range("G2").Value = 0.0060560181680545
$ws.Range("H2").Value = 0.978046934140802
$ws.Range("I2").Value = 0.0060560181680545
$ws.Range("J2").Value = 0.98107494322483
$ws.Range("K2").Value = 0.00227100681302044
$ws.Range("L2").Value = 0.00227100681302044
$ws.Range("M2").Value = 0.0408781226343679
$ws.Range("N2").Value = 0.00454201362604088
$ws.Range("O2").Value = 0.999242997728993
$ws.Range("P2").Value = 0.00529901589704769
$ws.Range("Q2").Value = 0.00151400454201363
$ws.Range("R2").Value = 0.00151400454201363
$ws.Range("S2").Value = 0.00151400454201363
$ws.Range("T2").Value = 0.0113550340651022
$ws.Range("U2").Value = 0.0060560181680545
$ws.Range("V2").Value = 0.0113550340651022
$ws.Range("W2").Value = 0.990915972747918
$ws.Range("X2").Value = 0.00302800908402725

$ws.Range("B3").Value = 0.0174110522331567
$ws.Range("C3").Value = 0.00151400454201363
$ws.Range("D3").Value = 0.00227100681302044
$ws.Range("E3").Value = 0.0158970476911431
$ws.Range("F3").Value = 0.98107494322483
$ws.Range("G3").Value = 0.978803936411809
$ws.Range("H3").Value = 0.00681302043906132
$ws.Range("I3").Value = 0.00757002271006813
$ws.Range("J3").Value = 0.00454201362604088
$ws.Range("K3").Value = 0.00302800908402725
$ws.Range("L3").Value = 0.851627554882665
$ws.Range("M3").Value = 0.0522331566994701
$ws.Range("N3").Value = 0.000757002271006813
$ws.Range("P3").Value = 0.0189250567751703
$ws.Range("Q3").Value = 0.985616956850871
$ws.Range("R3").Value = 0.995457986373959
$ws.Range("S3").Value = 0.0507191521574565
$ws.Range("U3").Value = 0.00378501135503407
$ws.Range("W3").Value = 0.000757002271006813
$ws.Range("X3").Value = 0.000757002271006813

$ws.Range("B4").Value = 0.0060560181680545
$ws.Range("C4").Value = 0.0287660862982589
$ws.Range("D4").Value = 0.973504920514762
$ws.Range("E4").Value = 0.00757002271006813
$ws.Range("G4").Value = 0.0136260408781226
$ws.Range("H4").Value = 0.00151400454201363
$ws.Range("I4").Value = 0.969719909159728
$ws.Range("J4").Value = 0.0136260408781226
$ws.Range("K4").Value = 0.00151400454201363
$ws.Range("L4").Value = 0.0158970476911431
$ws.Range("M4").Value = 0.00378501135503407
$ws.Range("N4").Value = 0.993943981831946
$ws.Range("O4").Value = 0.000757002271006813
$ws.Range("P4").Value = 0.974261922785768
$ws.Range("R4").Value = 0.00151400454201363
$ws.Range("S4").Value = 0.000757002271006813
$ws.Range("T4").Value = 0.987130961392884
$ws.Range("U4").Value = 0.989401968205905
$ws.Range("V4").Value = 0.976532929598789
$ws.Range("W4").Value = 0.00454201362604088
$ws.Range("X4").Value = 0.994700984102952

$ws.Range("B5").Value = 0.975018925056775
$ws.Range("C5").Value = 0.00454201362604088
$ws.Range("D5").Value = 0.0166540499621499
$ws.Range("E5").Value = 0.00151400454201363
$ws.Range("F5").Value = 0.0166540499621499
$ws.Range("G5").Value = 0.00151400454201363
$ws.Range("H5").Value = 0.0136260408781226
$ws.Range("I5").Value = 0.0166540499621499
$ws.Range("J5").Value = 0.000757002271006813
$ws.Range("K5").Value = 0.993186979560939
$ws.Range("L5").Value = 0.130204390613172
$ws.Range("M5").Value = 0.903103709311128
$ws.Range("N5").Value = 0.000757002271006813
$ws.Range("Q5").Value = 0.0128690386071158
$ws.Range("R5").Value = 0.00151400454201363
$ws.Range("S5").Value = 0.947009841029523
$ws.Range("T5").Value = 0.000757002271006813
$ws.Range("V5").Value = 0.012112036336109
$ws.Range("W5").Value = 0.00378501135503407
$ws.Range("X5").Value = 0.000757002271006813
